# lecturer.xlsx template update
#   - merge the "name" / "firstname" header columns into a single "_firstname"
#     column, reorder the remaining headers to _firstname, lastname, email, role
#   - shrink Table1 from 5 to 4 columns to match
#   - re-theme the sheet font from Calibri to Tahoma
#   - resize the header columns, refresh the saved view (zoom/selection)
#   - add an in-cell dropdown (data validation) listing the allowed roles

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-label the header row -------------------------------------------------
# Old layout: A=name  B=email  C=firstname  D=lastname  E=role
# New layout: A=_firstname  B=lastname  C=email  D=role  (E no longer used)
$ws.Range("A1").Value = "_firstname"
$ws.Range("B1").Value = "lastname"
$ws.Range("C1").Value = "email"
$ws.Range("D1").Value = "role"
$ws.Range("E1").ClearContents()

# --- 2. Shrink Table1 so it only covers the four remaining columns --------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D2"))

# --- 3. Swap the workbook's base font ------------------------------------------
$wb.Styles.Item("Normal").Font.Name = "Tahoma"

# --- 4. Column widths (character units) for the new layout ----------------------
$ws.Columns.Item(1).ColumnWidth = 14.75    # -> ~15.75 "_firstname"
$ws.Columns.Item(2).ColumnWidth = 13.25    # -> ~14.25 "lastname"
$ws.Columns.Item(3).ColumnWidth = 10.925   # -> ~11.875 "email"
$ws.Columns.Item(4).ColumnWidth = 14.255   # -> ~15.125 "role"

# --- 5. Refresh the saved view state ---------------------------------------------
$wb.Windows.Item(1).Zoom = 145
$ws.Range("I7").Select()

# --- 6. Restrict the role cell to a fixed dropdown list --------------------------
$ws.Range("D2").Validation.Add(3, 1, 1, '"lecturer, tabee manager, head of curriculum, moderator"')
